# Update "想去人数" (F column) counts across the four sheets to match the
# refreshed data snapshot (gh-pages output regenerated at 456a3b4).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 6743
$ws.Range("F4").Value = 1101
$ws.Range("F7").Value = 714
$ws.Range("F8").Value = 714
$ws.Range("F13").Value = 853
$ws.Range("F15").Value = 10
$ws.Range("F16").Value = 1013
$ws.Range("F17").Value = 1333
$ws.Range("F21").Value = 551
$ws.Range("F25").Value = 1058
$ws.Range("F26").Value = 1489
$ws.Range("F28").Value = 520
$ws.Range("F30").Value = 449
$ws.Range("F33").Value = 1127
$ws.Range("F37").Value = 1227
$ws.Range("F38").Value = 440
$ws.Range("F40").Value = 3851
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F8").Value = 25
$ws.Range("F11").Value = 156
$ws.Range("F25").Value = 231
$ws.Range("F27").Value = 108
$ws.Range("F29").Value = 227
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F5").Value = 1630
$ws.Range("F6").Value = 445
$ws.Range("F8").Value = 954
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 1630
$ws.Range("F5").Value = 445
$ws.Range("F7").Value = 954
$ws.Range("F9").Value = 6743
$ws.Range("F14").Value = 714
$ws.Range("F15").Value = 714
$ws.Range("F19").Value = 853
$ws.Range("F20").Value = 25
$ws.Range("F22").Value = 156
$ws.Range("F23").Value = 156
$ws.Range("F25").Value = 1013
$ws.Range("F26").Value = 1333
$ws.Range("F30").Value = 551
$ws.Range("F34").Value = 1058
$ws.Range("F35").Value = 1489
$ws.Range("F38").Value = 520
$ws.Range("F40").Value = 449
$ws.Range("F43").Value = 1127
$ws.Range("F46").Value = 227
$ws.Range("F49").Value = 1227
$ws.Range("F50").Value = 440
$ws.Range("F51").Value = 3851
